# Update loading_percent values for the 380 kV case (res_line data)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 7.729717733731942
$ws.Range("E2").Value = 12.72596180790705
$ws.Range("F2").Value = 44.58921583010526
$ws.Range("G2").Value = 57.84141350152662
$ws.Range("H2").Value = 19.8007200114177
$ws.Range("J2").Value = 10.10904146479808
$ws.Range("D3").Value = 7.731607383498007
$ws.Range("E3").Value = 12.73477324223081
$ws.Range("F3").Value = 43.7779840042736
$ws.Range("G3").Value = 56.12954056282797
$ws.Range("H3").Value = 19.59404963810029
$ws.Range("J3").Value = 10.11490316061086
$ws.Range("D4").Value = 7.732971663349006
$ws.Range("E4").Value = 12.74289888234465
$ws.Range("F4").Value = 43.28880581502413
$ws.Range("G4").Value = 55.07480788211523
$ws.Range("H4").Value = 19.47262906352413
$ws.Range("J4").Value = 10.12069478094531
$ws.Range("D5").Value = 7.733578766040466
$ws.Range("E5").Value = 12.74689087757522
$ws.Range("F5").Value = 43.09197246160716
$ws.Range("G5").Value = 54.64481053314546
$ws.Range("H5").Value = 19.42457083818628
$ws.Range("J5").Value = 10.12360443381651
$ws.Range("D6").Value = 7.733682659711601
$ws.Range("E6").Value = 12.74759479490289
$ws.Range("F6").Value = 43.05944780799565
$ws.Range("G6").Value = 54.57342058210156
$ws.Range("H6").Value = 19.41667790995177
$ws.Range("J6").Value = 10.1241207127091
$ws.Range("D7").Value = 7.732979644054461
$ws.Range("E7").Value = 12.74294996641732
$ws.Range("F7").Value = 43.28614073533018
$ws.Range("G7").Value = 55.0690085488706
$ws.Range("H7").Value = 19.47197512071685
$ws.Range("J7").Value = 10.12073179906159
$ws.Range("D8").Value = 7.730326842399804
$ws.Range("E8").Value = 12.7284352259098
$ws.Range("F8").Value = 44.30781285096185
$ws.Range("G8").Value = 57.25234060176189
$ws.Range("H8").Value = 19.72835136274261
$ws.Range("J8").Value = 10.1106064877137
$ws.Range("D9").Value = 7.726751176464111
$ws.Range("E9").Value = 12.72160251247393
$ws.Range("F9").Value = 46.37024980785323
$ws.Range("G9").Value = 61.47491567459098
$ws.Range("H9").Value = 20.272438758247
$ws.Range("J9").Value = 10.10822160243079
$ws.Range("D10").Value = 7.725127013220631
$ws.Range("E10").Value = 12.72987491738819
$ws.Range("F10").Value = 47.90568748517465
$ws.Range("G10").Value = 64.5041441229143
$ws.Range("H10").Value = 20.69436942026431
$ws.Range("J10").Value = 10.1172120709113
$ws.Range("D11").Value = 7.724608544374761
$ws.Range("E11").Value = 12.73654248767239
$ws.Range("F11").Value = 48.60533101872034
$ws.Range("G11").Value = 65.85944379355703
$ws.Range("H11").Value = 20.89045543231258
$ws.Range("J11").Value = 10.12365035974147
$ws.Range("D12").Value = 7.724444128017886
$ws.Range("E12").Value = 12.73948595990018
$ws.Range("F12").Value = 48.87017583043203
$ws.Range("G12").Value = 66.36888674787757
$ws.Range("H12").Value = 20.96524579847778
$ws.Range("J12").Value = 10.12642695060379
$ws.Range("D13").Value = 7.724478114804402
$ws.Range("E13").Value = 12.73883340025738
$ws.Range("F13").Value = 48.81314504986582
$ws.Range("G13").Value = 66.25934484517296
$ws.Range("H13").Value = 20.94911543145763
$ws.Range("J13").Value = 10.12581389385634
$ws.Range("D14").Value = 7.724594376744808
$ws.Range("E14").Value = 12.73677625479582
$ws.Range("F14").Value = 48.62712332969031
$ws.Range("G14").Value = 65.90143448747511
$ws.Range("H14").Value = 20.89659806435154
$ws.Range("J14").Value = 10.12387200243351
$ws.Range("D15").Value = 7.724669753947895
$ws.Range("E15").Value = 12.73557073299515
$ws.Range("F15").Value = 48.51315967612334
$ws.Range("G15").Value = 65.68169744160296
$ws.Range("H15").Value = 20.86449779800469
$ws.Range("J15").Value = 10.12272664843228
$ws.Range("D16").Value = 7.72516535002237
$ws.Range("E16").Value = 12.7294977078151
$ws.Range("F16").Value = 47.85996489177801
$ws.Range("G16").Value = 64.41507003753799
$ws.Range("H16").Value = 20.68163311876455
$ws.Range("J16").Value = 10.11683864795563
$ws.Range("D17").Value = 7.72552598314639
$ws.Range("E17").Value = 12.72651675488241
$ws.Range("F17").Value = 47.45935932342487
$ws.Range("G17").Value = 63.63183684729863
$ws.Range("H17").Value = 20.57047167034492
$ws.Range("J17").Value = 10.11382872043839
$ws.Range("D18").Value = 7.72575413793746
$ws.Range("E18").Value = 12.72507552557142
$ws.Range("F18").Value = 47.22906635946149
$ws.Range("G18").Value = 63.17923164279129
$ws.Range("H18").Value = 20.50692855796406
$ws.Range("J18").Value = 10.11231850326612
$ws.Range("D19").Value = 7.725834940557013
$ws.Range("E19").Value = 12.72463446245588
$ws.Range("F19").Value = 47.15112272176356
$ws.Range("G19").Value = 63.02564133938753
$ws.Range("H19").Value = 20.48548342974217
$ws.Range("J19").Value = 10.11184509972878
$ws.Range("D20").Value = 7.725485446234343
$ws.Range("E20").Value = 12.72680578544055
$ws.Range("F20").Value = 47.5019934939683
$ws.Range("G20").Value = 63.71543563327837
$ws.Range("H20").Value = 20.58226464743871
$ws.Range("J20").Value = 10.11412625151561
$ws.Range("D21").Value = 7.724559359727599
$ws.Range("E21").Value = 12.73736912101625
$ws.Range("F21").Value = 48.68176703450273
$ws.Range("G21").Value = 66.00666786147254
$ws.Range("H21").Value = 20.9120095951505
$ws.Range("J21").Value = 10.12443318921653
$ws.Range("D22").Value = 7.724140255688095
$ws.Range("E22").Value = 12.74671310399783
$ws.Range("F22").Value = 49.45214975722523
$ws.Range("G22").Value = 67.48188650198659
$ws.Range("H22").Value = 21.13061830241061
$ws.Range("J22").Value = 10.13314293819164
$ws.Range("D23").Value = 7.724346829107237
$ws.Range("E23").Value = 12.74150251249816
$ws.Range("F23").Value = 49.04112689914091
$ws.Range("G23").Value = 66.69672593583054
$ws.Range("H23").Value = 21.01367872353849
$ws.Range("J23").Value = 10.12831357741492
$ws.Range("D24").Value = 7.725503708133395
$ws.Range("E24").Value = 12.72667426583615
$ws.Range("F24").Value = 47.48271850152679
$ws.Range("G24").Value = 63.6776477902407
$ws.Range("H24").Value = 20.57693190082324
$ws.Range("J24").Value = 10.11399105181372
$ws.Range("D25").Value = 7.727543224493869
$ws.Range("E25").Value = 12.72112336097134
$ws.Range("F25").Value = 45.8076383838171
$ws.Range("G25").Value = 56.12954056282797
$ws.Range("H25").Value = 19.59404963810029
$ws.Range("J25").Value = 10.11490316061086
